# Add a new "NHSFLuVaccService" worksheet right after "LoginTest", populate
# it with the NHS Flu Vaccination service test data, and make it the active
# (selected) sheet/tab - matching the authored diff.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("LoginTest")

# Insert the new sheet right after LoginTest (so order becomes:
# LoginTest, NHSFLuVaccService, NewMedicineService, PatientRecords).
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$newSheet.Name = "NHSFLuVaccService"

# ---- Header row (bold, text-formatted) ----
$headerRange = $newSheet.Range("A1:F1")
$headerRange.NumberFormat = "@"
$headerRange.Font.Bold = $true

$newSheet.Range("A1").Value = "TestCaseName"
$newSheet.Range("B1").Value = "DateOfAdmin"
$newSheet.Range("E1").Value = "BatchNum"
$newSheet.Range("C1").Value = "TimeOfAdmin"
$newSheet.Range("D1").Value = "VaccineSearch"
$newSheet.Range("F1").Value = "ExpiryDate"

# ---- Data row (text-formatted) ----
$dataRange = $newSheet.Range("A2:F2")
$dataRange.NumberFormat = "@"

$newSheet.Range("A2").Value = "NHSFluVac1"
$newSheet.Range("B2").Value = "21/01/2021"
$newSheet.Range("C2").Value = "01:30 PM"
$newSheet.Range("D2").Value = "GSK - Fluarix TM Tetra (pack size 10)"
$newSheet.Range("E2").Value = "12345"
$newSheet.Range("F2").Value = "01/2021"

# ---- Column widths (closest achievable approximation of the authored
# widths: 16.08984375, 13.1796875, 14.08984375, 33.81640625, 11.453125,
# 10.7265625 character units) ----
$newSheet.Columns.Item(1).ColumnWidth = 15.3333333333
$newSheet.Columns.Item(2).ColumnWidth = 12.3333333333
$newSheet.Columns.Item(3).ColumnWidth = 13.3333333333
$newSheet.Columns.Item(4).ColumnWidth = 33.0
$newSheet.Columns.Item(5).ColumnWidth = 10.6666666667
$newSheet.Columns.Item(6).ColumnWidth = 9.8333333333

# Portrait page orientation, like the other sheets.
$newSheet.PageSetup.Orientation = 1

# Make this newly inserted sheet the active tab/selection, and select E17
# on it (matching the authored selection), which also clears
# tabSelected from the previously-active LoginTest sheet.
$newSheet.Activate()
$newSheet.Range("E17").Select() | Out-Null
